$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 8: continuation of Slider widget with another signal "scroll"
$ws.Range("C8").Value = "scroll"

# Add new row 9: new widget "radio buttons" with signal "checked" and slot
$ws.Range("B9").Value = "radio buttons"
$ws.Range("C9").Value = "checked"
$ws.Range("D9").Value = "switch between ROI selection modes"

# Update selection to match the new active cell D9
$ws.Range("D9").Select()
